$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.430374
$ws.Range("H2").Value = 7.291122000000001
$ws.Range("I2").Value = 0.009222757332915244
$ws.Range("J2").Value = 0.009222757332915246
$ws.Range("M2").Value = 20.94338233333333
$ws.Range("N2").Value = 62.830147
$ws.Range("O2").Value = 0.06014699761632732
$ws.Range("P2").Value = 0.06014699761632732
$ws.Range("Q2").Value = 50.90025189499266
$ws.Range("R2").Value = 458.102267054934
$ws.Range("S2").Value = 0.0005547211633188184
$ws.Range("T2").Value = 0.0005547211633188186
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.430374
$ws.Range("H3").Value = 7.291122000000001
$ws.Range("I3").Value = 0.009222757332915244
$ws.Range("J3").Value = 0.009222757332915246
$ws.Range("O3").Value = 0.105314686172636
$ws.Range("P3").Value = 0.105314686172636
$ws.Range("Q3").Value = 89.12405052407999
$ws.Range("R3").Value = 802.11645471672
$ws.Range("S3").Value = 0.0009712917941623464
$ws.Range("T3").Value = 0.0009712917941623468
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.430374
$ws.Range("H4").Value = 7.291122000000001
$ws.Range("I4").Value = 0.009222757332915244
$ws.Range("J4").Value = 0.009222757332915246
$ws.Range("M4").Value = 284.2332763333333
$ws.Range("N4").Value = 852.6998289999999
$ws.Range("O4").Value = 0.8162854462572833
$ws.Range("P4").Value = 0.8162854462572834
$ws.Range("Q4").Value = 690.7931647353487
$ws.Range("R4").Value = 6217.138482618137
$ws.Range("S4").Value = 0.007528402585221351
$ws.Range("T4").Value = 0.007528402585221354
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.430374
$ws.Range("H5").Value = 7.291122000000001
$ws.Range("I5").Value = 0.009222757332915244
$ws.Range("J5").Value = 0.009222757332915246
$ws.Range("M5").Value = 6.355709333333334
$ws.Range("N5").Value = 19.067128
$ws.Range("O5").Value = 0.01825286995375338
$ws.Range("P5").Value = 0.01825286995375338
$ws.Range("Q5").Value = 15.44675071529067
$ws.Range("R5").Value = 139.020756437616
$ws.Range("S5").Value = 0.0001683417902127273
$ws.Range("T5").Value = 0.0001683417902127274
$ws.Range("I6").Value = 0.5480399755605952
$ws.Range("J6").Value = 0.5480399755605954
$ws.Range("M6").Value = 20.94338233333333
$ws.Range("N6").Value = 62.830147
$ws.Range("O6").Value = 0.06014699761632732
$ws.Range("P6").Value = 0.06014699761632732
$ws.Range("Q6").Value = 3024.623959800361
$ws.Range("R6").Value = 27221.61563820325
$ws.Range("S6").Value = 0.03296295910369521
$ws.Range("T6").Value = 0.03296295910369521
$ws.Range("I7").Value = 0.5480399755605952
$ws.Range("J7").Value = 0.5480399755605954
$ws.Range("O7").Value = 0.105314686172636
$ws.Range("P7").Value = 0.105314686172636
$ws.Range("Q7").Value = 5295.980443588119
$ws.Range("R7").Value = 47663.82399229307
$ws.Range("S7").Value = 0.0577166580362232
$ws.Range("T7").Value = 0.05771665803622322
$ws.Range("I8").Value = 0.5480399755605952
$ws.Range("J8").Value = 0.5480399755605954
$ws.Range("M8").Value = 284.2332763333333
$ws.Range("N8").Value = 852.6998289999999
$ws.Range("O8").Value = 0.8162854462572833
$ws.Range("P8").Value = 0.8162854462572834
$ws.Range("Q8").Value = 41048.70761023479
$ws.Range("R8").Value = 369438.3684921131
$ws.Range("S8").Value = 0.4473570560173111
$ws.Range("T8").Value = 0.4473570560173112
$ws.Range("I9").Value = 0.5480399755605952
$ws.Range("J9").Value = 0.5480399755605954
$ws.Range("M9").Value = 6.355709333333334
$ws.Range("N9").Value = 19.067128
$ws.Range("O9").Value = 0.01825286995375338
$ws.Range("P9").Value = 0.01825286995375338
$ws.Range("Q9").Value = 917.885679837425
$ws.Range("R9").Value = 8260.971118536825
$ws.Range("S9").Value = 0.01000330240336572
$ws.Range("T9").Value = 0.01000330240336573
$ws.Range("G10").Value = 116.470388
$ws.Range("H10").Value = 349.411164
$ws.Range("I10").Value = 0.4419805861132828
$ws.Range("J10").Value = 0.4419805861132828
$ws.Range("M10").Value = 20.94338233333333
$ws.Range("N10").Value = 62.830147
$ws.Range("O10").Value = 0.06014699761632732
$ws.Range("P10").Value = 0.06014699761632732
$ws.Range("Q10").Value = 2439.283866395679
$ws.Range("R10").Value = 21953.5547975611
$ws.Range("S10").Value = 0.02658380525941857
$ws.Range("T10").Value = 0.02658380525941858
$ws.Range("G11").Value = 116.470388
$ws.Range("H11").Value = 349.411164
$ws.Range("I11").Value = 0.4419805861132828
$ws.Range("J11").Value = 0.4419805861132828
$ws.Range("O11").Value = 0.105314686172636
$ws.Range("P11").Value = 0.105314686172636
$ws.Range("Q11").Value = 4271.076280716959
$ws.Range("R11").Value = 38439.68652645264
$ws.Range("S11").Value = 0.04654704672091811
$ws.Range("T11").Value = 0.04654704672091812
$ws.Range("G12").Value = 116.470388
$ws.Range("H12").Value = 349.411164
$ws.Range("I12").Value = 0.4419805861132828
$ws.Range("J12").Value = 0.4419805861132828
$ws.Range("M12").Value = 284.2332763333333
$ws.Range("N12").Value = 852.6998289999999
$ws.Range("O12").Value = 0.8162854462572833
$ws.Range("P12").Value = 0.8162854462572834
$ws.Range("Q12").Value = 33104.75997705455
$ws.Range("R12").Value = 297942.8397934909
$ws.Range("S12").Value = 0.3607823199725367
$ws.Range("T12").Value = 0.3607823199725367
$ws.Range("G13").Value = 116.470388
$ws.Range("H13").Value = 349.411164
$ws.Range("I13").Value = 0.4419805861132828
$ws.Range("J13").Value = 0.4419805861132828
$ws.Range("M13").Value = 6.355709333333334
$ws.Range("N13").Value = 19.067128
$ws.Range("O13").Value = 0.01825286995375338
$ws.Range("P13").Value = 0.01825286995375338
$ws.Range("Q13").Value = 740.2519320685548
$ws.Range("R13").Value = 6662.267388616993
$ws.Range("S13").Value = 0.008067414160409448
$ws.Range("T13").Value = 0.00806741416040945
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1994
$ws.Range("H14").Value = 0.5982000000000001
$ws.Range("I14").Value = 0.0007566809932065188
$ws.Range("J14").Value = 0.0007566809932065189
$ws.Range("M14").Value = 20.94338233333333
$ws.Range("N14").Value = 62.830147
$ws.Range("O14").Value = 0.06014699761632732
$ws.Range("P14").Value = 0.06014699761632732
$ws.Range("Q14").Value = 4.176110437266667
$ws.Range("R14").Value = 37.5849939354
$ws.Range("S14").Value = 0.00004551208989471267
$ws.Range("T14").Value = 0.00004551208989471269
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1994
$ws.Range("H15").Value = 0.5982000000000001
$ws.Range("I15").Value = 0.0007566809932065188
$ws.Range("J15").Value = 0.0007566809932065189
$ws.Range("O15").Value = 0.105314686172636
$ws.Range("P15").Value = 0.105314686172636
$ws.Range("Q15").Value = 7.312181448
$ws.Range("R15").Value = 65.80963303199999
$ws.Range("S15").Value = 0.00007968962133234305
$ws.Range("T15").Value = 0.00007968962133234308
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1994
$ws.Range("H16").Value = 0.5982000000000001
$ws.Range("I16").Value = 0.0007566809932065188
$ws.Range("J16").Value = 0.0007566809932065189
$ws.Range("M16").Value = 284.2332763333333
$ws.Range("N16").Value = 852.6998289999999
$ws.Range("O16").Value = 0.8162854462572833
$ws.Range("P16").Value = 0.8162854462572834
$ws.Range("Q16").Value = 56.67611530086667
$ws.Range("R16").Value = 510.0850377078
$ws.Range("S16").Value = 0.0006176676822139876
$ws.Range("T16").Value = 0.0006176676822139877
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1994
$ws.Range("H17").Value = 0.5982000000000001
$ws.Range("I17").Value = 0.0007566809932065188
$ws.Range("J17").Value = 0.0007566809932065189
$ws.Range("M17").Value = 6.355709333333334
$ws.Range("N17").Value = 19.067128
$ws.Range("O17").Value = 0.01825286995375338
$ws.Range("P17").Value = 0.01825286995375338
$ws.Range("Q17").Value = 1.267328441066667
$ws.Range("R17").Value = 11.4059559696
$ws.Range("S17").Value = 0.00001381159976547553
$ws.Range("T17").Value = 0.00001381159976547554

Write-Output "Updated cells with new TPM-derived values."
